$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.991.04'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '1.745.53'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('D4').Value = '''0.9995'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''251.30'
$ws.Range('E5').Value = '  +7.95%  '
$ws.Range('D6').Value = '''0.9998'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').Value = '''0.5152'
$ws.Range('E7').Value = '  -2.31%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  +0.29%  '
$ws.Range('D10').Value = '1.746.65'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').Value = '''0.07239'
$ws.Range('E11').Value = '  +0.38%  '
$ws.Range('D12').Value = '''15.23'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '''0.6519'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('D14').Value = '''4.634'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').Value = '''77.91'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '''0.9998'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '''0.9993'
$ws.Range('E17').Value = '  -0.19%  '
$ws.Range('D18').Value = '26.014.69'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').Value = '  +2.74%  '
$ws.Range('D20').Value = '''0.000006809'
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = '1.966.98'
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('D22').Value = '''4.307'
$ws.Range('E22').Value = '  -0.51%  '
$ws.Range('D23').Value = '''8.686'
$ws.Range('E23').Value = '  -1.11%  '
$ws.Range('D24').Value = '''5.391'
$ws.Range('E24').Value = '  +3.83%  '
$ws.Range('D25').Value = '''135.97'
$ws.Range('E25').Value = '  -2.41%  '
$ws.Range('D26').Value = '''1.513'
$ws.Range('E26').Value = '  -0.77%  '
$ws.Range('D27').Value = '''15.30'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '''1.790'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').Value = '''106.11'
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('D30').Value = '''3.961'
$ws.Range('E30').Value = '  +5.50%  '
$ws.Range('D31').Value = '''0.08254'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('D32').Value = '''3.681'
$ws.Range('E32').Value = '  +1.01%  '
$ws.Range('D33').Value = '''0.04684'
$ws.Range('E33').Value = '  +3.36%  '
$ws.Range('D34').Value = '''2.655'
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('D35').Value = '''1.002'
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').Value = '''0.6276'
$ws.Range('E36').Value = '  -0.27%  '
$ws.Range('D37').Value = '''2.732'
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').Value = '''0.01614'
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('D39').Value = '''1.928'
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').Value = '''0.9993'
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').Value = '''100.69'
$ws.Range('E41').Value = '  +2.79%  '
$ws.Range('D42').Value = '''0.3896'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').Value = '''0.7606'
$ws.Range('E43').Value = '  +4.00%  '
$ws.Range('D44').Value = '''5.030'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '''6.360'
$ws.Range('E45').Value = '  +1.35%  '
$ws.Range('D46').Value = '''0.1133'
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').Value = '''55.60'
$ws.Range('E47').Value = '  +3.16%  '
$ws.Range('D48').Value = '''0.05232'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').Value = '''30.81'
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.569'
$ws.Range('E50').Value = '  -0.71%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').Value = '''0.3449'
$ws.Range('E51').Value = '  +0.04%  '
